# Insert a new column "Pan *" before the existing "Category *" column
# (which currently sits in column E), and populate it with PAN values
# for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift "Category *" / "Sub Category *" (and their data) one column to
# the right by inserting a new column at E.
$ws.Columns.Item(5).Insert()

# The inserted column should carry the same width as column D (the
# "Portfolio Company *" column it was split off from).
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# New header for the inserted column.
$ws.Range("E1").Value = "Pan *"

# New PAN data for each row (grouped by portfolio company, mirroring the
# existing Category values).
$ws.Range("E2").Value = "A11111111"
$ws.Range("E3").Value = "A11111111"
$ws.Range("E4").Value = "B11111111"
$ws.Range("E5").Value = "B11111111"

# Match the final selection left behind by the edit.
$ws.Range("E6").Select()
